$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the existing table down by two rows to make room for a new
#     "Validation of Models" title block at the top of the sheet. ---
$ws.Rows("1:2").Insert()

# --- New title row at the top: "Validation of Models" ---
$ws.Range("A1").Value = "Validation of Models"
$ws.Range("A1:F1").Font.Bold = $true
$ws.Range("A1:F1").HorizontalAlignment = -4108
$ws.Range("A1:F1").Merge()

# --- Highlight the lr=9E-3 / batch=5 result row (now row 23) in bold ---
$ws.Range("A23:F23").Font.Bold = $true

# --- New "Final Testing of Model" section appended under the table ---
$ws.Range("A31").Value = "Final Testing of Model"
$ws.Range("A31:F31").Font.Bold = $true
$ws.Range("A31:F31").HorizontalAlignment = -4108
$ws.Range("A31:F31").Merge()

$ws.Range("A32").Value = 5
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = 0.0009
$ws.Range("D32").Value = 86.69
$ws.Range("E32").Value = 84.7
$ws.Range("F32").Value = 88.61
$ws.Range("A32:F32").HorizontalAlignment = -4108

$ws.Range("A33").Value = 5
$ws.Range("B33").Value = 5
$ws.Range("C33").Value = 0.0009
$ws.Range("D33").Value = 98.73
$ws.Range("E33").Value = 98.1
$ws.Range("F33").Value = 99.27
$ws.Range("A33:F33").HorizontalAlignment = -4108

# --- Column width tweaks ---
$ws.Columns("A").ColumnWidth = 10.666666666666666
$ws.Columns("C").ColumnWidth = 22.333333333333332
$ws.Columns("D").ColumnWidth = 22.166666666666668
$ws.Columns("E").ColumnWidth = 20.333333333333332
$ws.Columns("F").ColumnWidth = 22.333333333333332

# --- Print orientation ---
$ws.PageSetup.Orientation = 1

# --- Active cell / selection ---
$ws.Range("K34").Select()
